$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the crypto price refresh diff.
# D-column (Price) values look numeric (e.g. "63.529.38", "1.00", "0.0₃0859")
# so we force text storage via NumberFormat "@" before assigning, then restore
# the default "Normal" style so no stray style index is left on the cell.
# E-column (Volume/1h) values already carry padding spaces, which keeps Excel
# from reinterpreting them as numbers, so they are assigned directly.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.529.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.584.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.049.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.383.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("E16").Value = "  +5.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.570.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "554.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0859"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "166.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("E43").Value = "  +5.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0593"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("E46").Value = "  +6.29%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0963"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0233"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.05%  "
